$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, pushing existing rows 33:151 down to 34:152
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with this week's record
$ws.Cells.Item(33, 1).Value = 10
$ws.Cells.Item(33, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(33, 3).Value = "La Araucanía"
$ws.Cells.Item(33, 4).Value = [DateTime]"2022-12-15"
$ws.Cells.Item(33, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(33, 5).Value = 9
$ws.Cells.Item(33, 6).Value = 100114002
$ws.Cells.Item(33, 7).Value = "Camote"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 50
$ws.Cells.Item(33, 11).Value = 24000
$ws.Cells.Item(33, 12).Value = 24000
$ws.Cells.Item(33, 13).Value = 24000
$ws.Cells.Item(33, 14).Value = "`$/malla 20 kilos"
$ws.Cells.Item(33, 15).Value = "Perú"
$ws.Cells.Item(33, 16).Value = 1200
$ws.Cells.Item(33, 17).Value = 20
$ws.Cells.Item(33, 18).Value = "Hortaliza"
